# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.096.92"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "1.665.10"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").Value = "'209.51"
$ws.Range("E5").Value = "  -4.02%  "
$ws.Range("D6").Value = "'0.5168"
$ws.Range("E6").Value = "  -4.47%  "
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").Value = "'0.2626"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").Value = "'0.06197"
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("D10").Value = "'21.01"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("D11").Value = "'0.07503"
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("D12").Value = "1.661.54"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").Value = "'4.397"
$ws.Range("E13").Value = "  -2.85%  "
$ws.Range("D14").Value = "'0.5557"
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").Value = "'0.000007866"
$ws.Range("E15").Value = "  -5.76%  "
$ws.Range("D16").Value = "'65.41"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "26.108.38"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "'4.768"
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("D20").Value = "'10.34"
$ws.Range("E20").Value = "  -5.74%  "
$ws.Range("D21").Value = "'185.52"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("D22").Value = "'6.119"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "'147.04"
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("D25").Value = "'0.1234"
$ws.Range("E25").Value = "  -5.53%  "
$ws.Range("D26").Value = "'7.522"
$ws.Range("E26").Value = "  -4.16%  "
$ws.Range("D27").Value = "'15.78"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "'0.06220"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("D29").Value = "'1.355"
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("D30").Value = "'1.272"
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("D31").Value = "'3.468"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").Value = "'3.401"
$ws.Range("E32").Value = "  -5.00%  "
$ws.Range("D33").Value = "'1.610"
$ws.Range("E33").Value = "  -4.29%  "
$ws.Range("D34").Value = "'0.9903"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").Value = "'0.6007"
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").Value = "'2.703"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'6.101"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "1.070.92"
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("D41").Value = "'0.8583"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "  -1.21%  "
$ws.Range("D43").Value = "'98.81"
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("D44").Value = "1.812.47"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "'55.78"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").Value = "'0.05249"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").Value = "'7.933"
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("D50").Value = "'0.4251"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").Value = "'5.882"
$ws.Range("E51").Value = "  -2.43%  "
